$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '43.865.78'
$ws.Range('E2').Value = '  -5.54%  '

# Row 3
$ws.Range('D3').Value = '2.604.00'
$ws.Range('E3').Value = '  +0.24%  '

# Row 4
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '301.77'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.83%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '97.01'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.31%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.579'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -4.07%  '

# Row 8
$ws.Range('E8').Value = '  +0.14%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.561'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -3.28%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '37.32'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -5.13%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0817'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -3.49%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.88'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -3.81%  '

# Row 13
$ws.Range('D13').Value = '3.003.29'
$ws.Range('E13').Value = '  +0.29%  '

# Row 14
$ws.Range('E14').Value = '  +1.32%  '

# Row 15
$ws.Range('D15').Value = '2.602.31'
$ws.Range('E15').Value = '  -0.16%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.896'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -3.14%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.46'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -3.73%  '

# Row 18
$ws.Range('D18').Value = '43.819.26'
$ws.Range('E18').Value = '  -5.97%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.68'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.97%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0981'
$ws.Range('E20').Value = '  -3.53%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.47'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -4.07%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '73.32'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.46%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '267.26'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.82%  '

# Row 24
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.23'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.24%  '

# Row 25
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.95'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -3.20%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '29.69'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.40%  '

# Row 27
$ws.Range('E27').Value = '  +0.20%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.31'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -3.12%  '

# Row 29
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.22'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.03%  '

# Row 30
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '38.00'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.67%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.06'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -5.04%  '

# Row 32
$ws.Range('E32').Value = '  +0.15%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.27'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.33%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '152.07'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.45%  '

# Row 35
$ws.Range('E35').Value = '  -1.50%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0817'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.71%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.118'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -5.00%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '24.49'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +5.64%  '

# Row 39
$ws.Range('E39').Value = '  -1.25%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '16.92'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +5.43%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.55'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -2.15%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0316'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -4.39%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.88'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -4.79%  '

# Row 44
$ws.Range('D44').Value = '2.075.26'
$ws.Range('E44').Value = '  -2.68%  '

# Row 45
$ws.Range('E45').Value = '  -0.36%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '88.69'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.64%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.22'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.94%  '

# Row 48
$ws.Range('E48').Value = '  +3.71%  '

# Row 49
$ws.Range('D49').Value = '2.853.78'
$ws.Range('E49').Value = '  +0.26%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '106.63'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.43%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.192'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -4.49%  '
